$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): values look numeric (contain dots) so Excel would
# auto-convert them to numbers on assignment. Force the cells to Text first,
# write the literal strings, then clear the format back to the default so the
# cells keep matching the original (unstyled) appearance, just like the
# source workbook where these are plain inline strings.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.437.71"
$ws.Range("D3").Value = "1.858.38"
$ws.Range("D5").Value = "311.08"
$ws.Range("D6").Value = "1.009"
$ws.Range("D7").Value = "0.4767"
$ws.Range("D8").Value = "0.3797"
$ws.Range("D9").Value = "0.07300"
$ws.Range("D10").Value = "0.9289"
$ws.Range("D11").Value = "20.71"
$ws.Range("D12").Value = "0.07783"
$ws.Range("D13").Value = "1.866.78"
$ws.Range("D14").Value = "5.437"
$ws.Range("D15").Value = "6.539"
$ws.Range("D16").Value = "90.12"
$ws.Range("D17").Value = "1.011"
$ws.Range("D18").Value = "0.000008804"
$ws.Range("D19").Value = "1.008"
$ws.Range("D20").Value = "27.461.39"
$ws.Range("D21").Value = "14.61"
$ws.Range("D22").Value = "5.090"
$ws.Range("D23").Value = "10.67"
$ws.Range("D24").Value = "1.942"
$ws.Range("D25").Value = "154.77"
$ws.Range("D26").Value = "18.43"
$ws.Range("D27").Value = "1.997"
$ws.Range("D28").Value = "115.37"
$ws.Range("D29").Value = "4.940"
$ws.Range("D30").Value = "0.08884"
$ws.Range("D31").Value = "3.327"
$ws.Range("D32").Value = "1.201"
$ws.Range("D33").Value = "0.7522"
$ws.Range("D34").Value = "4.578"
$ws.Range("D35").Value = "2.692"
$ws.Range("D37").Value = "0.02036"
$ws.Range("D38").Value = "0.5535"
$ws.Range("D39").Value = "0.05269"
$ws.Range("D40").Value = "2.986"
$ws.Range("D41").Value = "7.013"
$ws.Range("D42").Value = "8.565"
$ws.Range("D43").Value = "0.1512"
$ws.Range("D44").Value = "0.4857"
$ws.Range("D45").Value = "10.63"
$ws.Range("D46").Value = "1.010"
$ws.Range("D47").Value = "1.661"
$ws.Range("D48").Value = "102.90"
$ws.Range("D49").Value = "67.30"
$ws.Range("D50").Value = "0.06088"
$ws.Range("D51").Value = "0.9130"

$dRange.ClearFormats()

# --- Volume(1h) column (E): values already contain "%" and padding spaces,
# so Excel keeps them as text without any extra nudging.
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E37").Value = "  +3.92%  "
$ws.Range("E38").Value = "  +5.60%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +2.80%  "

